# Update the semivariogram-models worksheet: add GO, MG and MT states,
# update MS and PA figures, and re-order the data rows accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows, in final order (row 2 .. row 8)
# columns: state, year, model, c0, c0_c1, a, gde, rss, r2
$rows = @(
    @("GO", 2015, "Sph", 0,      2.1275, 0.6,                0,                   2878.8526,          0.6132),
    @("MG", 2015, "Exp", 0,      1.5897, 0.6899999999999999, 0,                   2876.9595,          0.6907),
    @("MS", 2015, "Exp", 0.3214, 1.2175, 0.82,                0.2639835728952772, 621.2714999999999,  0.6555),
    @("MS", 2016, "Sph", 0.0793, 0.9015, 0.5600000000000001,  0.08796450360510261,212.1107,           0.8929),
    @("MT", 2015, "Exp", 0,      0.9879, 0.73,                0,                   2601.9701,          0.2146),
    @("PA", 2015, "Sph", 0.1198, 1.0045, 0.43,                0.1192633150821304, 420.9166,           0.3708),
    @("PA", 2016, "Sph", 0.7436, 1.4251, 0.23,                0.5217879447056347, 209752.6425,        0.856)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
    $ws.Cells.Item($r, 8).Value = $vals[7]
    $ws.Cells.Item($r, 9).Value = $vals[8]
}
